# Confidence intervals for geom_abline
# Adds "min"/"max" header + per-row confidence-interval bounds (columns E:G)
# to the summary table at the bottom of Sheet1 (rows 153-163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the min/max columns
$ws.Range("E153").Value = "min"
$ws.Range("F153").Value = "max"

# Min / max confidence-interval bounds for each hex-pair summary row
$data = @(
  @(154, 1,   12),
  @(155, 14,  25),
  @(156, 27,  38),
  @(157, 40,  56),
  @(158, 60,  74),
  @(159, 76,  88),
  @(160, 90,  102),
  @(161, 107, 122),
  @(162, 124, 141),
  @(163, 147, 162)
)

foreach ($entry in $data) {
  $r = $entry[0]
  $minVal = $entry[1]
  $maxVal = $entry[2]

  $ws.Range("E$r").NumberFormat = "0"
  $ws.Range("F$r").NumberFormat = "0"
  $ws.Range("G$r").NumberFormat = "0"

  $ws.Range("E$r").Value = $minVal
  $ws.Range("F$r").Value = $maxVal
}

# Update the frozen-pane view / active selection to match the new extent
$ws.Range("A163").Select()
